$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update DAMSLTag (column I) and DialogAct (column J) values for re-run SGNN dialog act annotations
$ws.Range("I2").Value = 'b'
$ws.Range("J2").Value = 'Acknowledge (Backchannel)'
$ws.Range("I9").Value = 'sd'
$ws.Range("J9").Value = 'Statement-non-opinion'
$ws.Range("I11").Value = '%'
$ws.Range("J11").Value = 'Uninterpretable'
$ws.Range("I22").Value = 'aa'
$ws.Range("J22").Value = 'Agree/Accept'
$ws.Range("I28").Value = 'sd'
$ws.Range("J28").Value = 'Statement-non-opinion'
$ws.Range("I30").Value = 'sd'
$ws.Range("J30").Value = 'Statement-non-opinion'
$ws.Range("I42").Value = 'sd'
$ws.Range("J42").Value = 'Statement-non-opinion'
$ws.Range("I46").Value = 'b'
$ws.Range("J46").Value = 'Acknowledge (Backchannel)'
$ws.Range("I48").Value = 'b'
$ws.Range("J48").Value = 'Acknowledge (Backchannel)'
$ws.Range("I55").Value = 'sd'
$ws.Range("J55").Value = 'Statement-non-opinion'
$ws.Range("I59").Value = 'sd'
$ws.Range("J59").Value = 'Statement-non-opinion'
$ws.Range("I61").Value = 'sv'
$ws.Range("J61").Value = 'Statement-opinion'
$ws.Range("I71").Value = 'sv'
$ws.Range("J71").Value = 'Statement-opinion'
$ws.Range("I79").Value = 'sv'
$ws.Range("J79").Value = 'Statement-opinion'
$ws.Range("I83").Value = 'sd'
$ws.Range("J83").Value = 'Statement-non-opinion'
$ws.Range("I91").Value = 'sv'
$ws.Range("J91").Value = 'Statement-opinion'
$ws.Range("I117").Value = 'aa'
$ws.Range("J117").Value = 'Agree/Accept'
$ws.Range("I121").Value = 'sv'
$ws.Range("J121").Value = 'Statement-opinion'
$ws.Range("I123").Value = 'aa'
$ws.Range("J123").Value = 'Agree/Accept'
$ws.Range("I152").Value = 'sd'
$ws.Range("J152").Value = 'Statement-non-opinion'
$ws.Range("I170").Value = 'sd'
$ws.Range("J170").Value = 'Statement-non-opinion'
$ws.Range("I171").Value = 'sv'
$ws.Range("J171").Value = 'Statement-opinion'
$ws.Range("I193").Value = 'sv'
$ws.Range("J193").Value = 'Statement-opinion'
$ws.Range("I206").Value = 'sd'
$ws.Range("J206").Value = 'Statement-non-opinion'
$ws.Range("I223").Value = 'sv'
$ws.Range("J223").Value = 'Statement-opinion'
$ws.Range("I228").Value = 'sd'
$ws.Range("J228").Value = 'Statement-non-opinion'
$ws.Range("I232").Value = 'sd'
$ws.Range("J232").Value = 'Statement-non-opinion'
$ws.Range("I235").Value = 'sv'
$ws.Range("J235").Value = 'Statement-opinion'
$ws.Range("I238").Value = 'sv'
$ws.Range("J238").Value = 'Statement-opinion'
$ws.Range("I239").Value = 'sd'
$ws.Range("J239").Value = 'Statement-non-opinion'
$ws.Range("I242").Value = 'sd'
$ws.Range("J242").Value = 'Statement-non-opinion'
$ws.Range("I263").Value = 'ba'
$ws.Range("J263").Value = 'Appreciation'
$ws.Range("I279").Value = 'sd'
$ws.Range("J279").Value = 'Statement-non-opinion'
$ws.Range("I284").Value = 'sd'
$ws.Range("J284").Value = 'Statement-non-opinion'
$ws.Range("I290").Value = 'sv'
$ws.Range("J290").Value = 'Statement-opinion'
$ws.Range("I305").Value = 'sv'
$ws.Range("J305").Value = 'Statement-opinion'
$ws.Range("I308").Value = 'sv'
$ws.Range("J308").Value = 'Statement-opinion'
$ws.Range("I311").Value = 'sv'
$ws.Range("J311").Value = 'Statement-opinion'
$ws.Range("I313").Value = 'sv'
$ws.Range("J313").Value = 'Statement-opinion'
$ws.Range("I316").Value = 'sd'
$ws.Range("J316").Value = 'Statement-non-opinion'
$ws.Range("I319").Value = 'sd'
$ws.Range("J319").Value = 'Statement-non-opinion'
$ws.Range("I329").Value = 'sd'
$ws.Range("J329").Value = 'Statement-non-opinion'
$ws.Range("I345").Value = 'sd'
$ws.Range("J345").Value = 'Statement-non-opinion'
$ws.Range("I346").Value = 'sd'
$ws.Range("J346").Value = 'Statement-non-opinion'
$ws.Range("I347").Value = 'sd'
$ws.Range("J347").Value = 'Statement-non-opinion'
$ws.Range("I366").Value = 'b'
$ws.Range("J366").Value = 'Acknowledge (Backchannel)'
$ws.Range("I373").Value = 'sd'
$ws.Range("J373").Value = 'Statement-opinion'
